# Update changelog.xlsx to reflect the GLSL 460 shader version changes:
# add a new changelog entry row documenting the soul fire model fix (Sodium).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new changelog bullet in column A, row 13 (right after the
# existing "1.21.7" row), inheriting the column's default style.
$ws.Range("A13").Value = "Fixed soul fire model (Sodium)"

# Move/record the active selection to the next empty row, as in the saved file.
$ws.Range("A14").Select()
